$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Datos actualizados" timestamp footer in A1
$ws.Range("A1").Value = "Datos actualizados a 30 de Abril de 2020 a las 18:52"

# Update country rows (re-sorted order + refreshed case counts)
$ws.Range("A4").Value = "Estados Unidos"
$ws.Range("B4").Value = 1073588
$ws.Range("C4").Value = 9394
$ws.Range("D4").Value = 149604
$ws.Range("E4").Value = 861721
$ws.Range("F4").Value = 18697
$ws.Range("G4").Value = 608
$ws.Range("H4").Value = 62263

$ws.Range("A15").Value = "Canada"
$ws.Range("B15").Value = 52069
$ws.Range("C15").Value = 472
$ws.Range("D15").Value = 20936
$ws.Range("E15").Value = 28051
$ws.Range("F15").Value = 557
$ws.Range("G15").Value = 86
$ws.Range("H15").Value = 3082

$ws.Range("A55").Value = "Marruecos"
$ws.Range("B55").Value = 4423
$ws.Range("C55").Value = 102
$ws.Range("D55").Value = 984
$ws.Range("E55").Value = 3269
$ws.Range("F55").Value = 1
$ws.Range("G55").Value = 2
$ws.Range("H55").Value = 170

$ws.Range("A60").Value = "Luxemburgo"
$ws.Range("B60").Value = 3784
$ws.Range("C60").Value = 15
$ws.Range("D60").Value = 3213
$ws.Range("E60").Value = 481
$ws.Range("F60").Value = 23
$ws.Range("G60").Value = 1
$ws.Range("H60").Value = 90

$ws.Range("A139").Value = "Gibraltar"
$ws.Range("B139").Value = 144
$ws.Range("C139").Value = 3
$ws.Range("D139").Value = 131
$ws.Range("E139").Value = 13
$ws.Range("F139").Value = 0
$ws.Range("G139").Value = 0
$ws.Range("H139").Value = 0

$ws.Range("A140").Value = "Liberia"
$ws.Range("B140").Value = 141
$ws.Range("C140").Value = 0
$ws.Range("D140").Value = 45
$ws.Range("E140").Value = 80
$ws.Range("F140").Value = 0
$ws.Range("G140").Value = 0
$ws.Range("H140").Value = 16

$ws.Range("A151").Value = "Zambia"
$ws.Range("B151").Value = 106
$ws.Range("C151").Value = 9
$ws.Range("D151").Value = 55
$ws.Range("E151").Value = 48
$ws.Range("F151").Value = 1
$ws.Range("G151").Value = 0
$ws.Range("H151").Value = 3

$ws.Range("A152").Value = "Suazilandia"
$ws.Range("B152").Value = 100
$ws.Range("C152").Value = 9
$ws.Range("D152").Value = 12
$ws.Range("E152").Value = 87
$ws.Range("F152").Value = 0
$ws.Range("G152").Value = 0
$ws.Range("H152").Value = 1

$ws.Range("A153").Value = "Aruba"
$ws.Range("B153").Value = 100
$ws.Range("C153").Value = 0
$ws.Range("D153").Value = 73
$ws.Range("E153").Value = 25
$ws.Range("F153").Value = 4
$ws.Range("G153").Value = 0
$ws.Range("H153").Value = 2

$ws.Range("A154").Value = "Monaco"
$ws.Range("B154").Value = 95
$ws.Range("C154").Value = 0
$ws.Range("D154").Value = 64
$ws.Range("E154").Value = 27
$ws.Range("F154").Value = 1
$ws.Range("G154").Value = 0
$ws.Range("H154").Value = 4

$ws.Range("A201").Value = "Santa Sede"
$ws.Range("B201").Value = 11
$ws.Range("C201").Value = 1
$ws.Range("D201").Value = 2
$ws.Range("E201").Value = 9
$ws.Range("F201").Value = 0
$ws.Range("G201").Value = 0
$ws.Range("H201").Value = 0

$ws.Range("A202").Value = "Montserrat"
$ws.Range("B202").Value = 11
$ws.Range("C202").Value = 0
$ws.Range("D202").Value = 2
$ws.Range("E202").Value = 8
$ws.Range("F202").Value = 1
$ws.Range("G202").Value = 0
$ws.Range("H202").Value = 1

$ws.Range("A203").Value = "Burundi"
$ws.Range("B203").Value = 11
$ws.Range("C203").Value = 0
$ws.Range("D203").Value = 4
$ws.Range("E203").Value = 6
$ws.Range("F203").Value = 0
$ws.Range("G203").Value = 0
$ws.Range("H203").Value = 1

$ws.Range("A204").Value = "Seychelles"
$ws.Range("B204").Value = 11
$ws.Range("C204").Value = 0
$ws.Range("D204").Value = 6
$ws.Range("E204").Value = 5
$ws.Range("F204").Value = 0
$ws.Range("G204").Value = 0
$ws.Range("H204").Value = 0

$ws.Range("A205").Value = "Gambia"
$ws.Range("B205").Value = 11
$ws.Range("C205").Value = 1
$ws.Range("D205").Value = 8
$ws.Range("E205").Value = 2
$ws.Range("F205").Value = 0
$ws.Range("G205").Value = 0
$ws.Range("H205").Value = 1

$ws.Range("A206").Value = "Groenlandia"
$ws.Range("B206").Value = 11
$ws.Range("C206").Value = 0
$ws.Range("D206").Value = 11
$ws.Range("E206").Value = 0
$ws.Range("F206").Value = 0
$ws.Range("G206").Value = 0
$ws.Range("H206").Value = 0
